# AR-70-38.xlsx — "added a unit scaling function"
#
# The rain-drop model on the "rain" sheet blends two exponential
# distributions using a mixing weight in $I$5. The author tightened that
# weight from 0.99515546875000571 to the rounder 0.9951, which ripples
# through every dependent formula in the sheet (rows 4, 7-71). They also
# added a small "1/mean" scaling helper (mirroring the existing G5/H5
# inverse-scale pair) in H14:L14, and left the selection sitting on G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rain")

# Core input change: tighten the mixing-weight constant. Every formula that
# depends on it (B4:G4, G7:G12, F7:F71, I7:I9, etc.) recalculates from this.
$ws.Range("I5").Value = 0.9951

# New "unit scaling" helper block, mirroring the existing 1/mean columns
# (G5 -> I label/value pair, H5 -> L label/value pair). Write the label
# cells in K14 before H14 so the shared-string table gets "mean1" (idx 45)
# ahead of "mean" (idx 46), matching the source order.
$ws.Range("K14").Value = "mean1"
$ws.Range("H14").Value = "mean"
$ws.Range("I14").Formula = "=1/G5"
$ws.Range("L14").Formula = "=1/H5"

# Match the formatting of the existing analogous cells (G5 -> H14:I14,
# H5 -> K14:L14) instead of inventing new styles.
$ws.Range("G5").Copy()
$ws.Range("H14:I14").PasteSpecial(-4122)

$ws.Range("H5").Copy()
$ws.Range("K14:L14").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Leave the sheet selection on G11, as in the saved file.
$ws.Activate()
[void]$ws.Range("G11").Select()
